$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of row 6 into row 7 so the new row matches existing style (s="3"/s="4")
$ws.Range("A6:B6").Copy() | Out-Null
$ws.Range("A7:B7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Set the new row's values
$ws.Range("A7").Value = "28-09-2025"
$ws.Range("B7").Value = $ws.Range("B6").Value2

$excel.CutCopyMode = 0
